# Applies the crypto price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.128.22'
$ws.Range('D3').Value = '3.475.58'
$ws.Range('E3').Value = '  -1.25%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''587.41'
$ws.Range('E5').Value = '  -3.22%  '
$ws.Range('D6').Value = '''136.67'
$ws.Range('E6').Value = '  -4.87%  '
$ws.Range('D7').Value = '3.474.17'
$ws.Range('E7').Value = '  -1.30%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  -3.50%  '
$ws.Range('E10').Value = '  -6.20%  '
$ws.Range('D11').Value = '''7.13'
$ws.Range('E11').Value = '  -7.28%  '
$ws.Range('E12').Value = '  -5.98%  '
$ws.Range('D13').Value = '4.066.16'
$ws.Range('E13').Value = '  -1.10%  '
$ws.Range('E14').Value = '  -6.71%  '
$ws.Range('D15').Value = '3.488.30'
$ws.Range('E15').Value = '  -0.99%  '
$ws.Range('D16').Value = '''26.43'
$ws.Range('E16').Value = '  -7.83%  '
$ws.Range('E17').Value = '  -1.35%  '
$ws.Range('D18').Value = '65.108.51'
$ws.Range('E18').Value = '  -1.96%  '
$ws.Range('E19').Value = '  -9.68%  '
$ws.Range('E20').Value = '  -5.91%  '
$ws.Range('D21').Value = '''13.89'
$ws.Range('E21').Value = '  -4.91%  '
$ws.Range('D22').Value = '''388.15'
$ws.Range('E22').Value = '  -8.30%  '
$ws.Range('D23').Value = '''0.554'
$ws.Range('E23').Value = '  -6.02%  '
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').Value = '''72.47'
$ws.Range('E25').Value = '  -5.86%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').Value = '3.616.18'
$ws.Range('E27').Value = '  -1.39%  '
$ws.Range('E28').Value = '  -4.38%  '
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('D30').Value = '''7.32'
$ws.Range('E30').Value = '  -6.58%  '
$ws.Range('D31').Value = '''8.18'
$ws.Range('E31').Value = '  -8.01%  '
$ws.Range('E32').Value = '  -10.14%  '
$ws.Range('D33').Value = '3.494.02'
$ws.Range('E33').Value = '  -0.93%  '
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('E35').Value = '  -6.06%  '
$ws.Range('D36').Value = '''23.08'
$ws.Range('E36').Value = '  -4.58%  '
$ws.Range('D37').Value = '''171.41'
$ws.Range('E37').Value = '  -1.23%  '
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D38').Value = '''6.81'
$ws.Range('E38').Value = '  -9.63%  '
$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D39').Value = '''1.19'
$ws.Range('E39').Value = '  -10.62%  '
$ws.Range('E40').Value = '  -10.04%  '
$ws.Range('D41').Value = '''4.74'
$ws.Range('E41').Value = '  -8.70%  '
$ws.Range('D42').Value = '''0.0778'
$ws.Range('E42').Value = '  -3.85%  '
$ws.Range('D43').Value = '''0.811'
$ws.Range('E43').Value = '  -4.67%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').Value = '''1.00'
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').Value = '''42.45'
$ws.Range('E45').Value = '  -6.66%  '
$ws.Range('D46').Value = '''24.68'
$ws.Range('E46').Value = '  +7.32%  '
$ws.Range('D47').Value = '''4.35'
$ws.Range('E47').Value = '  -12.63%  '
$ws.Range('E48').Value = '  -9.15%  '
$ws.Range('E49').Value = '  +2.93%  '
$ws.Range('E50').Value = '  -5.47%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '2.220.05'
$ws.Range('E51').Value = '  -3.75%  '

Write-Host "Applied cryptos update: 89 cells"
